$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells: "_old" -> "_FV2404" and "_new" -> "_FV2410"
$headers = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headers.Count; $i++) {
    # columns A..J (1..10) -> "_old" suffix becomes "_FV2404"
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i] + "_FV2404"
    # columns L..U (12..21) -> "_new" suffix becomes "_FV2410"
    $ws.Cells.Item(1, $i + 12).Value = $headers[$i] + "_FV2410"
}
# column K (11) stays "diff"
$ws.Cells.Item(1, 11).Value = "diff"

# 2) Turn the data range into an Excel Table ("Table1") over A1:U64
$tableRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $tableRange, $false, 1)
$table.Name = "Table1"

# 3) Freeze the header row (top row)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
